$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$rows1 = @{
    2  = 2328
    3  = 471
    5  = 331
    6  = 331
    7  = 536
    9  = 757
    11 = 765
    13 = 86
    15 = 15
    16 = 1018
    17 = 19685
    18 = 576
    19 = 59
    20 = 216
    21 = 281
    22 = 168
    23 = 141
    24 = 12
    25 = 6
    26 = 201
    28 = 324
    29 = 137
}
foreach ($row in $rows1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $rows1[$row]
}

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$rows2 = @{
    7  = 220
    8  = 3382
    10 = 86
    14 = 120
    16 = 3013
}
foreach ($row in $rows2.Keys) {
    $ws2.Cells.Item($row, 6).Value = $rows2[$row]
}

# Sheet "本地生活" (local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$rows3 = @{
    2 = 265
    3 = 88
    4 = 558
    5 = 202
}
foreach ($row in $rows3.Keys) {
    $ws3.Cells.Item($row, 6).Value = $rows3[$row]
}

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$rows4 = @{
    2  = 265
    3  = 88
    5  = 2328
    6  = 558
    7  = 471
    9  = 331
    10 = 331
    11 = 536
    17 = 202
    18 = 757
    20 = 765
    22 = 86
    24 = 15
    25 = 1018
    26 = 19685
    27 = 220
    28 = 3382
    30 = 86
    32 = 576
    33 = 59
    34 = 216
    37 = 281
    38 = 168
    39 = 141
    40 = 12
    41 = 6
    42 = 120
    44 = 201
    46 = 324
    47 = 137
    48 = 3013
}
foreach ($row in $rows4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $rows4[$row]
}
